$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(201).Insert()

$ws.Cells.Item(201, 1).Value = 3
$ws.Cells.Item(201, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(201, 3).Value = "Coquimbo"
$ws.Cells.Item(201, 4).Value = 44907
$ws.Cells.Item(201, 5).Value = 5
$ws.Cells.Item(201, 6).Value = 100112039
$ws.Cells.Item(201, 7).Value = "Ciboulette"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 120
$ws.Cells.Item(201, 11).Value = 1500
$ws.Cells.Item(201, 12).Value = 1500
$ws.Cells.Item(201, 13).Value = 1500
$ws.Cells.Item(201, 14).Value = "$/docena de atados"
$ws.Cells.Item(201, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(201, 16).Value = 500
$ws.Cells.Item(201, 17).Value = 3
$ws.Cells.Item(201, 18).Value = "Hortaliza"
